$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row before row 29; this shifts old rows 29-35
# (and their formatting) down to rows 30-36.
$ws.Rows("29:29").Insert()

# Grow the table (表3) to cover the newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:D36"))

# The "UseDiamond" column is text-formatted ("false"/"true" stored as
# shared strings, not booleans). Copy the cell that was pushed down to
# row 30 (still holding the "false" text cell) into the new row 29 so it
# keeps the same text representation instead of being auto-typed as a
# boolean.
$ws.Range("D30").Copy($ws.Range("D29"))

# Re-sequence the Id column (A) so it stays a contiguous run
# 15000036 .. 15000043 regardless of which Item ended up on which row.
$ws.Range("A29").Value = 15000036
$ws.Range("A30").Value = 15000037
$ws.Range("A31").Value = 15000038
$ws.Range("A32").Value = 15000039
$ws.Range("A33").Value = 15000040
$ws.Range("A34").Value = 15000041
$ws.Range("A35").Value = 15000042
$ws.Range("A36").Value = 15000043

# Fill in the new row's other cells (Item / Shelf).
$ws.Range("B29").Value = "jueyu"
$ws.Range("C29").Value = 3

# Mirror the author's final cursor position / scroll state in the saved view.
$null = $ws.Range("B26").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
